$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column ("Price") values are plain numeric-looking text (e.g. "1.01",
# "27.513.72") that must stay literal strings, matching the source sheet's
# inline-string cells. Excel auto-coerces Range.Value assignments that look
# like numbers into real numbers, so we force the cell to Text format first
# and restore the default "Normal" style afterwards so no stray number
# formatting / style index is left behind on the cell.

$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = "27.513.72"
$r.Style = "Normal"
$ws.Range("E2").Value = "  -1.40%  "
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = "1.574.01"
$r.Style = "Normal"
$ws.Range("E3").Value = "  -3.53%  "
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = "1.01"
$r.Style = "Normal"
$ws.Range("E4").Value = "  +0.52%  "
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = "205.99"
$r.Style = "Normal"
$ws.Range("E5").Value = "  -2.55%  "
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = "0.500"
$r.Style = "Normal"
$ws.Range("E6").Value = "  -3.54%  "
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = "1.01"
$r.Style = "Normal"
$ws.Range("E7").Value = "  +0.61%  "
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = "22.10"
$r.Style = "Normal"
$ws.Range("E8").Value = "  -5.99%  "
$ws.Range("E9").Value = "  -2.29%  "
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = "0.0588"
$r.Style = "Normal"
$ws.Range("E10").Value = "  -3.98%  "
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = "0.0865"
$r.Style = "Normal"
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = "1.806.01"
$r.Style = "Normal"
$ws.Range("E12").Value = "  -3.05%  "
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = "1.570.39"
$r.Style = "Normal"
$ws.Range("E13").Value = "  -3.73%  "
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = "3.83"
$r.Style = "Normal"
$ws.Range("E14").Value = "  -4.99%  "
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = "0.530"
$r.Style = "Normal"
$ws.Range("E15").Value = "  -5.87%  "
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = "27.503.31"
$r.Style = "Normal"
$ws.Range("E16").Value = "  -1.44%  "
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = "62.73"
$r.Style = "Normal"
$ws.Range("E17").Value = "  -4.02%  "
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = "215.84"
$r.Style = "Normal"
$ws.Range("E18").Value = "  -5.77%  "
$r = $ws.Range("D19")
$r.NumberFormat = "@"
$r.Value = "7.29"
$r.Style = "Normal"
$ws.Range("E19").Value = "  -4.81%  "
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = "0.0₃0688"
$r.Style = "Normal"
$ws.Range("E20").Value = "  -4.41%  "
$ws.Range("E21").Value = "  +0.51%  "
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = "4.13"
$r.Style = "Normal"
$ws.Range("E22").Value = "  -4.75%  "
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = "9.49"
$r.Style = "Normal"
$ws.Range("E23").Value = "  -5.69%  "
$ws.Range("E24").Value = "  -4.06%  "
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = "153.26"
$r.Style = "Normal"
$ws.Range("E25").Value = "  -0.94%  "
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = "1.01"
$r.Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "
$r = $ws.Range("D27")
$r.NumberFormat = "@"
$r.Value = "6.70"
$r.Style = "Normal"
$ws.Range("E27").Value = "  -2.68%  "
$r = $ws.Range("D28")
$r.NumberFormat = "@"
$r.Value = "14.97"
$r.Style = "Normal"
$ws.Range("E28").Value = "  -3.57%  "
$ws.Range("E29").Value = "  -4.62%  "
$r = $ws.Range("D30")
$r.NumberFormat = "@"
$r.Value = "1.15"
$r.Style = "Normal"
$ws.Range("E30").Value = "  -2.63%  "
$r = $ws.Range("D31")
$r.NumberFormat = "@"
$r.Value = "0.0464"
$r.Style = "Normal"
$ws.Range("E31").Value = "  -3.58%  "
$r = $ws.Range("D32")
$r.NumberFormat = "@"
$r.Value = "3.23"
$r.Style = "Normal"
$ws.Range("E32").Value = "  -5.39%  "
$r = $ws.Range("D33")
$r.NumberFormat = "@"
$r.Value = "1.362.40"
$r.Style = "Normal"
$ws.Range("E33").Value = "  -2.33%  "
$ws.Range("E34").Value = "  -5.48%  "
$ws.Range("E35").Value = "  -5.56%  "
$r = $ws.Range("D36")
$r.NumberFormat = "@"
$r.Value = "0.963"
$r.Style = "Normal"
$ws.Range("E36").Value = "  -5.75%  "
$ws.Range("E37").Value = "  -1.21%  "
$r = $ws.Range("D38")
$r.NumberFormat = "@"
$r.Value = "0.0164"
$r.Style = "Normal"
$ws.Range("E38").Value = "  -3.68%  "
$ws.Range("E39").Value = "  -4.23%  "
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = "0.806"
$r.Style = "Normal"
$ws.Range("E40").Value = "  -4.90%  "
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = "1.01"
$r.Style = "Normal"
$ws.Range("E41").Value = "  +0.55%  "
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = "0.974"
$r.Style = "Normal"
$ws.Range("E42").Value = "  -3.89%  "
$r = $ws.Range("D43")
$r.NumberFormat = "@"
$r.Value = "2.22"
$r.Style = "Normal"
$ws.Range("E43").Value = "  +3.42%  "
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = "1.76"
$r.Style = "Normal"
$ws.Range("E44").Value = "  -3.88%  "
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = "5.25"
$r.Style = "Normal"
$ws.Range("E45").Value = "  -3.37%  "
$r = $ws.Range("D46")
$r.NumberFormat = "@"
$r.Value = "63.10"
$r.Style = "Normal"
$ws.Range("E46").Value = "  -4.21%  "
$r = $ws.Range("D47")
$r.NumberFormat = "@"
$r.Value = "1.713.48"
$r.Style = "Normal"
$ws.Range("E47").Value = "  -3.34%  "
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = "86.71"
$r.Style = "Normal"
$ws.Range("E48").Value = "  -2.25%  "
$r = $ws.Range("D49")
$r.NumberFormat = "@"
$r.Value = "0.0965"
$r.Style = "Normal"
$ws.Range("E49").Value = "  -4.88%  "
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("E51").Value = "  -7.04%  "
